$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 102.4929963333333
$ws.Range("H2").Value = 307.478989
$ws.Range("I2").Value = 0.2065071987599813
$ws.Range("J2").Value = 0.2065071987599814
$ws.Range("M2").Value = 2.724001666666667
$ws.Range("N2").Value = 8.172005
$ws.Range("O2").Value = 0.04635500474236593
$ws.Range("P2").Value = 0.04635500474236593
$ws.Range("Q2").Value = 279.1910928336605
$ws.Range("R2").Value = 2512.719835502945
$ws.Range("S2").Value = 0.009572642177851639
$ws.Range("T2").Value = 0.009572642177851639
$ws.Range("G3").Value = 102.4929963333333
$ws.Range("H3").Value = 307.478989
$ws.Range("I3").Value = 0.2065071987599813
$ws.Range("J3").Value = 0.2065071987599814
$ws.Range("O3").Value = 0.6912512390256352
$ws.Range("P3").Value = 0.6912512390256351
$ws.Range("Q3").Value = 4163.330149976355
$ws.Range("R3").Value = 37469.97134978719
$ws.Range("S3").Value = 0.1427483570105502
$ws.Range("T3").Value = 0.1427483570105502
$ws.Range("G4").Value = 102.4929963333333
$ws.Range("H4").Value = 307.478989
$ws.Range("I4").Value = 0.2065071987599813
$ws.Range("J4").Value = 0.2065071987599814
$ws.Range("M4").Value = 15.419285
$ws.Range("N4").Value = 46.257855
$ws.Range("O4").Value = 0.2623937562319988
$ws.Range("P4").Value = 0.2623937562319988
$ws.Range("Q4").Value = 1580.368720967622
$ws.Range("R4").Value = 14223.3184887086
$ws.Range("S4").Value = 0.05418619957157948
$ws.Range("T4").Value = 0.05418619957157948
$ws.Range("I5").Value = 0.581825957350084
$ws.Range("J5").Value = 0.5818259573500841
$ws.Range("M5").Value = 2.724001666666667
$ws.Range("N5").Value = 8.172005
$ws.Range("O5").Value = 0.04635500474236593
$ws.Range("P5").Value = 0.04635500474236593
$ws.Range("Q5").Value = 786.6099867073488
$ws.Range("R5").Value = 7079.489880366139
$ws.Range("S5").Value = 0.02697054501219474
$ws.Range("T5").Value = 0.02697054501219475
$ws.Range("I6").Value = 0.581825957350084
$ws.Range("J6").Value = 0.5818259573500841
$ws.Range("O6").Value = 0.6912512390256352
$ws.Range("P6").Value = 0.6912512390256351
$ws.Range("S6").Value = 0.402187913915522
$ws.Range("T6").Value = 0.402187913915522
$ws.Range("I7").Value = 0.581825957350084
$ws.Range("J7").Value = 0.5818259573500841
$ws.Range("M7").Value = 15.419285
$ws.Range("N7").Value = 46.257855
$ws.Range("O7").Value = 0.2623937562319988
$ws.Range("P7").Value = 0.2623937562319988
$ws.Range("Q7").Value = 4452.627073363326
$ws.Range("R7").Value = 40073.64366026993
$ws.Range("S7").Value = 0.1526674984223673
$ws.Range("T7").Value = 0.1526674984223673
$ws.Range("G8").Value = 105.053815
$ws.Range("H8").Value = 315.161445
$ws.Range("I8").Value = 0.2116668438899346
$ws.Range("J8").Value = 0.2116668438899346
$ws.Range("M8").Value = 2.724001666666667
$ws.Range("N8").Value = 8.172005
$ws.Range("O8").Value = 0.04635500474236593
$ws.Range("P8").Value = 0.04635500474236593
$ws.Range("Q8").Value = 286.1667671496917
$ws.Range("R8").Value = 2575.500904347225
$ws.Range("S8").Value = 0.009811817552319548
$ws.Range("T8").Value = 0.00981181755231955
$ws.Range("G9").Value = 105.053815
$ws.Range("H9").Value = 315.161445
$ws.Range("I9").Value = 0.2116668438899346
$ws.Range("J9").Value = 0.2116668438899346
$ws.Range("O9").Value = 0.6912512390256352
$ws.Range("P9").Value = 0.6912512390256351
$ws.Range("Q9").Value = 4267.352219239327
$ws.Range("R9").Value = 38406.16997315395
$ws.Range("S9").Value = 0.146314968099563
$ws.Range("T9").Value = 0.146314968099563
$ws.Range("G10").Value = 105.053815
$ws.Range("H10").Value = 315.161445
$ws.Range("I10").Value = 0.2116668438899346
$ws.Range("J10").Value = 0.2116668438899346
$ws.Range("M10").Value = 15.419285
$ws.Range("N10").Value = 46.257855
$ws.Range("O10").Value = 0.2623937562319988
$ws.Range("P10").Value = 0.2623937562319988
$ws.Range("Q10").Value = 1619.854713822275
$ws.Range("R10").Value = 14578.69242440047
$ws.Range("S10").Value = 0.05554005823805206
$ws.Range("T10").Value = 0.05554005823805207
